# Add Darwin plot & add Darwin to Combined
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) All_SANs - append two new asset rows (Desktop Mini G9 SANs)
# ---------------------------------------------------------------------
$wsAllSANs = $wb.Worksheets.Item("All_SANs")
$wsAllSANs.Cells.Item(119, 1).Value = "SAN125998"
$wsAllSANs.Cells.Item(119, 2).Value = "Desktop Mini G9"
$wsAllSANs.Cells.Item(119, 3).Value = "2024-11-17 13:16:43"

$wsAllSANs.Cells.Item(120, 1).Value = "SAN125999"
$wsAllSANs.Cells.Item(120, 2).Value = "Desktop Mini G9"
$wsAllSANs.Cells.Item(120, 3).Value = "2024-11-17 13:16:47"

# the sheet's column A:C has a default style (centered) applied at the
# column level; the newly appended rows use the plain/no-style format
# like the rows immediately above them, so reset it explicitly.
$wsAllSANs.Range("A119:C120").Style = "Normal"

# ---------------------------------------------------------------------
# 2) 4.2_Items - update Desktop Mini G9 last/new count
# ---------------------------------------------------------------------
$ws42Items = $wb.Worksheets.Item("4.2_Items")
$ws42Items.Cells.Item(2, 2).Value = 56
$ws42Items.Cells.Item(2, 3).Value = 58

# ---------------------------------------------------------------------
# 3) 4.2_Timestamps - append the two Desktop Mini G9 "add" log rows
#    plus a follow-up "add 2" row
# ---------------------------------------------------------------------
$ws42Timestamps = $wb.Worksheets.Item("4.2_Timestamps")
$ws42Timestamps.Cells.Item(28, 1).Value = "2024-11-17 13:16:43"
$ws42Timestamps.Cells.Item(28, 2).Value = "Desktop Mini G9"
$ws42Timestamps.Cells.Item(28, 3).Value = "add"
$ws42Timestamps.Cells.Item(28, 4).Value = "SAN125998"

$ws42Timestamps.Cells.Item(29, 1).Value = "2024-11-17 13:16:47"
$ws42Timestamps.Cells.Item(29, 2).Value = "Desktop Mini G9"
$ws42Timestamps.Cells.Item(29, 3).Value = "add"
$ws42Timestamps.Cells.Item(29, 4).Value = "SAN125999"

$ws42Timestamps.Cells.Item(30, 1).Value = "2024-11-17 13:16:47"
$ws42Timestamps.Cells.Item(30, 2).Value = "Desktop Mini G9"
$ws42Timestamps.Cells.Item(30, 3).Value = "add 2"

# ---------------------------------------------------------------------
# 4) BR_Items - Laptop Charger count
# ---------------------------------------------------------------------
$wsBRItems = $wb.Worksheets.Item("BR_Items")
$wsBRItems.Cells.Item(9, 2).Value = 5
$wsBRItems.Cells.Item(9, 3).Value = 35

# ---------------------------------------------------------------------
# 5) BR_Timestamps - append Laptop Charger add log rows
# ---------------------------------------------------------------------
$wsBRTimestamps = $wb.Worksheets.Item("BR_Timestamps")
$wsBRTimestamps.Cells.Item(20, 1).Value = "2024-11-17 13:19:45"
$wsBRTimestamps.Cells.Item(20, 2).Value = "Laptop Charger" + [char]160
$wsBRTimestamps.Cells.Item(20, 3).Value = "add 5"

$wsBRTimestamps.Cells.Item(21, 1).Value = "2024-11-17 13:19:57"
$wsBRTimestamps.Cells.Item(21, 2).Value = "Laptop Charger" + [char]160
$wsBRTimestamps.Cells.Item(21, 3).Value = "add 30"

# ---------------------------------------------------------------------
# 6) L17_Items - seed counts for Laptop 840 G6, Monitor 24, Monitor 34
# ---------------------------------------------------------------------
$wsL17Items = $wb.Worksheets.Item("L17_Items")
$wsL17Items.Cells.Item(2, 2).Value = 0
$wsL17Items.Cells.Item(2, 3).Value = 15

$wsL17Items.Cells.Item(3, 2).Value = 0
$wsL17Items.Cells.Item(3, 3).Value = 3

$wsL17Items.Cells.Item(4, 2).Value = 0
$wsL17Items.Cells.Item(4, 3).Value = 30

# ---------------------------------------------------------------------
# 7) L17_Timestamps - append the initial add log rows
# ---------------------------------------------------------------------
$wsL17Timestamps = $wb.Worksheets.Item("L17_Timestamps")
$wsL17Timestamps.Cells.Item(2, 1).Value = "2024-11-17 12:22:50"
$wsL17Timestamps.Cells.Item(2, 2).Value = "Monitor 24" + [char]8221 + [char]160
$wsL17Timestamps.Cells.Item(2, 3).Value = "add 3"

$wsL17Timestamps.Cells.Item(3, 1).Value = "2024-11-17 13:20:29"
$wsL17Timestamps.Cells.Item(3, 2).Value = "Monitor 34" + [char]8221 + " Ultrawide"
$wsL17Timestamps.Cells.Item(3, 3).Value = "add 30"

$wsL17Timestamps.Cells.Item(4, 1).Value = "2024-11-17 13:20:35"
$wsL17Timestamps.Cells.Item(4, 2).Value = "Laptop 840 G6"
$wsL17Timestamps.Cells.Item(4, 3).Value = "add 15"

# ---------------------------------------------------------------------
# 9) B4.3_Timestamps - clear the stale empty "SAN #" placeholder cell
#    (no value, Excel drops it on save)
# ---------------------------------------------------------------------
$wsB43Timestamps = $wb.Worksheets.Item("B4.3_Timestamps")
$wsB43Timestamps.Cells.Item(3, 4).Value = ""

# ---------------------------------------------------------------------
# 10) Darwin_Items - seed counts for the new site
# ---------------------------------------------------------------------
$wsDarwinItems = $wb.Worksheets.Item("Darwin_Items")
$wsDarwinItems.Cells.Item(9, 2).Value = 0
$wsDarwinItems.Cells.Item(9, 3).Value = 30

$wsDarwinItems.Cells.Item(11, 2).Value = 0
$wsDarwinItems.Cells.Item(11, 3).Value = 30

$wsDarwinItems.Cells.Item(12, 2).Value = 0
$wsDarwinItems.Cells.Item(12, 3).Value = 30

$wsDarwinItems.Cells.Item(13, 2).Value = 0
$wsDarwinItems.Cells.Item(13, 3).Value = 30

$wsDarwinItems.Cells.Item(15, 2).Value = 0
$wsDarwinItems.Cells.Item(15, 3).Value = 30

# ---------------------------------------------------------------------
# 11) Darwin_Timestamps - append the initial add log rows for Darwin
# ---------------------------------------------------------------------
$wsDarwinTimestamps = $wb.Worksheets.Item("Darwin_Timestamps")

# clear the stale empty "SAN #" placeholder cells left over on the
# existing rows (they carried no value, Excel drops them on save)
$wsDarwinTimestamps.Cells.Item(2, 4).Value = ""
$wsDarwinTimestamps.Cells.Item(3, 4).Value = ""
$wsDarwinTimestamps.Cells.Item(4, 4).Value = ""

$wsDarwinTimestamps.Cells.Item(5, 1).Value = "2024-11-17 13:20:05"
$wsDarwinTimestamps.Cells.Item(5, 2).Value = "Laptop Charger" + [char]160
$wsDarwinTimestamps.Cells.Item(5, 3).Value = "add 30"

$wsDarwinTimestamps.Cells.Item(6, 1).Value = "2024-11-17 13:20:08"
$wsDarwinTimestamps.Cells.Item(6, 2).Value = "Monitor 24" + [char]8221 + [char]160
$wsDarwinTimestamps.Cells.Item(6, 3).Value = "add 30"

$wsDarwinTimestamps.Cells.Item(7, 1).Value = "2024-11-17 13:20:09"
$wsDarwinTimestamps.Cells.Item(7, 2).Value = "Monitor 34" + [char]8221 + " Ultrawide"
$wsDarwinTimestamps.Cells.Item(7, 3).Value = "add 30"

$wsDarwinTimestamps.Cells.Item(8, 1).Value = "2024-11-17 13:20:11"
$wsDarwinTimestamps.Cells.Item(8, 2).Value = "USB DVD-RW Drive"
$wsDarwinTimestamps.Cells.Item(8, 3).Value = "add 30"

$wsDarwinTimestamps.Cells.Item(9, 1).Value = "2024-11-17 13:20:12"
$wsDarwinTimestamps.Cells.Item(9, 2).Value = "Wired Keyboard"
$wsDarwinTimestamps.Cells.Item(9, 3).Value = "add 30"
